$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$ws.Range("D2").Value = "'245.19"
$ws.Range("D2").Style = "Normal"
$ws.Range("E2").Value = "'-0.58%"
$ws.Range("E2").Style = "Normal"
$ws.Range("D3").Value = "'29.18"
$ws.Range("D3").Style = "Normal"
$ws.Range("E3").Value = "'-2.11%"
$ws.Range("E3").Style = "Normal"
$ws.Range("D4").Value = "'5.251"
$ws.Range("D4").Style = "Normal"
$ws.Range("E4").Value = "'1.67%"
$ws.Range("E4").Style = "Normal"
$ws.Range("D5").Value = "'0.05701"
$ws.Range("D5").Style = "Normal"
$ws.Range("E5").Value = "'-0.22%"
$ws.Range("E5").Style = "Normal"
$ws.Range("D6").Value = "'6.616"
$ws.Range("D6").Style = "Normal"
$ws.Range("E6").Value = "'0.03%"
$ws.Range("E6").Style = "Normal"
$ws.Range("D7").Value = "'3.192"
$ws.Range("D7").Style = "Normal"
$ws.Range("E7").Value = "'3.83%"
$ws.Range("E7").Style = "Normal"
$ws.Range("D8").Value = "'0.8511"
$ws.Range("D8").Style = "Normal"
$ws.Range("E8").Value = "'-0.58%"
$ws.Range("E8").Style = "Normal"
$ws.Range("D9").Value = "'0.8532"
$ws.Range("D9").Style = "Normal"
$ws.Range("E9").Value = "'-1.77%"
$ws.Range("E9").Style = "Normal"
$ws.Range("D10").Value = "'0.1370"
$ws.Range("D10").Style = "Normal"
$ws.Range("E10").Value = "'0.38%"
$ws.Range("E10").Style = "Normal"
$ws.Range("D11").Value = "'0.07069"
$ws.Range("D11").Style = "Normal"
$ws.Range("E11").Value = "'-0.11%"
$ws.Range("E11").Style = "Normal"
$ws.Range("B12").Value = "'LiechtensteinCryptoassetsExchange"
$ws.Range("B12").Style = "Normal"
$ws.Range("C12").Value = "'https://coinranking.com/coin/v4IW9oaF+liechtensteincryptoassetsexchange-lcx"
$ws.Range("C12").Style = "Normal"
$ws.Range("D12").Value = "'0.03242"
$ws.Range("D12").Style = "Normal"
$ws.Range("E12").Value = "'-3.83%"
$ws.Range("E12").Style = "Normal"
$ws.Range("B13").Value = "'BitrueCoin"
$ws.Range("B13").Style = "Normal"
$ws.Range("C13").Value = "'https://coinranking.com/coin/SLYjzF4ty+bitruecoin-btr"
$ws.Range("C13").Style = "Normal"
$ws.Range("D13").Value = "'0.03191"
$ws.Range("D13").Style = "Normal"
$ws.Range("E13").Value = "'9.24%"
$ws.Range("E13").Style = "Normal"
$ws.Range("B14").Value = "'BitMartToken"
$ws.Range("B14").Style = "Normal"
$ws.Range("C14").Value = "'https://coinranking.com/coin/6uzcPMFgWUJNH+bitmarttoken-bmx"
$ws.Range("C14").Style = "Normal"
$ws.Range("D14").Value = "'0.09234"
$ws.Range("D14").Style = "Normal"
$ws.Range("E14").Value = "'-1.62%"
$ws.Range("E14").Style = "Normal"
$ws.Range("B15").Value = "'BitForexToken"
$ws.Range("B15").Style = "Normal"
$ws.Range("C15").Value = "'https://coinranking.com/coin/2nh5ugplNocUp+bitforextoken-bf"
$ws.Range("C15").Style = "Normal"
$ws.Range("D15").Value = "'0.001523"
$ws.Range("D15").Style = "Normal"
$ws.Range("E15").Value = "'0.75%"
$ws.Range("E15").Style = "Normal"
$ws.Range("B16").Value = "'One"
$ws.Range("B16").Style = "Normal"
$ws.Range("C16").Value = "'https://coinranking.com/coin/6Lga5NiXX3rT+one-one"
$ws.Range("C16").Style = "Normal"
$ws.Range("D16").Value = "'0.0005963"
$ws.Range("D16").Style = "Normal"
$ws.Range("E16").Value = "'-0.62%"
$ws.Range("E16").Style = "Normal"
$ws.Range("B17").Value = "'TigerCash"
$ws.Range("B17").Style = "Normal"
$ws.Range("C17").Value = "'https://coinranking.com/coin/6hIn06L2+tigercash-tch"
$ws.Range("C17").Style = "Normal"
$ws.Range("D17").Value = "'0.005992"
$ws.Range("D17").Style = "Normal"
$ws.Range("E17").Value = "'-2.59%"
$ws.Range("E17").Style = "Normal"
$ws.Range("B18").Value = "'LEO"
$ws.Range("B18").Style = "Normal"
$ws.Range("C18").Value = "'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo"
$ws.Range("C18").Style = "Normal"
$ws.Range("D18").Value = "'3.495"
$ws.Range("D18").Style = "Normal"
$ws.Range("E18").Value = "'0.24%"
$ws.Range("E18").Style = "Normal"
$ws.Range("B19").Value = "'BTSEToken"
$ws.Range("B19").Style = "Normal"
$ws.Range("C19").Value = "'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse"
$ws.Range("C19").Style = "Normal"
$ws.Range("D19").Value = "'2.175"
$ws.Range("D19").Style = "Normal"
$ws.Range("E19").Value = "'-0.64%"
$ws.Range("E19").Style = "Normal"
$ws.Range("B20").Value = "'BitpandaEcosystemToken"
$ws.Range("B20").Style = "Normal"
$ws.Range("C20").Value = "'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best"
$ws.Range("C20").Style = "Normal"
$ws.Range("D20").Value = "'0.3158"
$ws.Range("D20").Style = "Normal"
$ws.Range("E20").Value = "'-0.50%"
$ws.Range("E20").Style = "Normal"
$ws.Range("E21").Value = "'-1.79%"
$ws.Range("E21").Style = "Normal"
$ws.Range("D22").Value = "'3.493"
$ws.Range("D22").Style = "Normal"
$ws.Range("E22").Value = "'0.70%"
$ws.Range("E22").Style = "Normal"
$ws.Range("D23").Value = "'0.04084"
$ws.Range("D23").Style = "Normal"
$ws.Range("E23").Value = "'-2.25%"
$ws.Range("E23").Style = "Normal"
$ws.Range("E24").Value = "'-0.01%"
$ws.Range("E24").Style = "Normal"
$ws.Range("D25").Value = "'0.001223"
$ws.Range("D25").Style = "Normal"
$ws.Range("E25").Value = "'0.24%"
$ws.Range("E25").Style = "Normal"
$ws.Range("D26").Value = "'0.004141"
$ws.Range("D26").Style = "Normal"
$ws.Range("E26").Value = "'-17.62%"
$ws.Range("E26").Style = "Normal"
$ws.Range("D27").Value = "'0.0001200"
$ws.Range("D27").Style = "Normal"
$ws.Range("E27").Value = "'-0.80%"
$ws.Range("E27").Style = "Normal"
$ws.Range("E28").Value = "'-25.22%"
$ws.Range("E28").Style = "Normal"
$ws.Range("D40").Value = "'0.03753"
$ws.Range("D40").Style = "Normal"
$ws.Range("E40").Value = "'0.19%"
$ws.Range("E40").Style = "Normal"
$ws.Range("D41").Value = "'0.1066"
$ws.Range("D41").Style = "Normal"
$ws.Range("E41").Value = "'-0.67%"
$ws.Range("E41").Style = "Normal"
$ws.Range("D42").Value = "'0.003712"
$ws.Range("D42").Style = "Normal"
$ws.Range("E42").Value = "'-35.56%"
$ws.Range("E42").Style = "Normal"
$ws.Range("D43").Value = "'0.002402"
$ws.Range("D43").Style = "Normal"
$ws.Range("D44").Value = "'0.009363"
$ws.Range("D44").Style = "Normal"
$ws.Range("E44").Value = "'2.01%"
$ws.Range("E44").Style = "Normal"
$ws.Range("E45").Value = "'1.52%"
$ws.Range("E45").Style = "Normal"
$ws.Range("D47").Value = "'0.07504"
$ws.Range("D47").Style = "Normal"
